# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that get refreshed each time the
# handback status report is (re)generated.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for first row (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-20 15:06:42"

# "zh-cn" sheet: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-20 15:06:38"
$wsZhCn.Range("K2").Value = "2016-08-20 15:06:56"

# "de-de" sheet: Correspond Handoff Datetime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-20 15:07:05"
